# Upates to require angle brackets for @base and @prefix values.
#
# The "Info" worksheet lists @base/@prefix declarations; the URI/URL values
# in column D need to be wrapped in angle brackets (<...>). Also bump the
# column D width so the longer values are visible, and move the sheet
# selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Wrap the four URI/URL values in column D with angle brackets. Written in
# this particular order so new shared-string entries land in the same
# sequence as the canonical workbook (ex1, model#, data#, short-data#).
$ws.Range("D1").Value2 = "<http://example.org/ex1>"
$ws.Range("D3").Value2 = "<http://foo.bar/model#>"
$ws.Range("D2").Value2 = "<http://foo.bar/data#>"
$ws.Range("D4").Value2 = "<http:/foo.bar/data#>"
$ws.Range("D5").Value2 = "<http://foo.bar/model#>"

# Widen column D to comfortably fit the longer, bracketed values (closest
# the host's pixel-quantised ColumnWidth can get to the canonical
# 38.5703125 OOXML width is 38.5, reached via any input in ~[37.58,37.75]).
$ws.Columns.Item(4).ColumnWidth = 37.67

# Move the active selection to D5.
$ws.Activate()
$ws.Range("D5").Select()
